# Move the 5 "New" listings into "Previously added" (appended at the end),
# then clear the "New" sheet back down to just its header row.

$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# Hyperlink targets for the 5 rows currently on "New" (A2:A6), in row order.
$urls = @(
    "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/valgundes-nov/cghpfg.html",
    "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/livberzes-pag/ohbhg.html",
    "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/valgundes-nov/egcjx.html",
    "https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/valgundes-nov/ecmgi.html",
    "https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/varmes-pag/dpmoh.html"
)

# 1) Copy the 5 data rows (values + formatting) onto the end of "Previously added".
$destStartRow = $wsPrev.UsedRange.Rows.Count() + 1
$lastRow = $destStartRow + $urls.Count - 1
$srcRange = $wsNew.Range("A2:F6")
$destRange = $wsPrev.Range("A" + $destStartRow)
$srcRange.Copy($destRange)

# 2) Recreate the hyperlinks on column A of the newly appended rows.
for ($i = 0; $i -lt $urls.Count; $i++) {
    $row = $destStartRow + $i
    $cell = $wsPrev.Range("A" + $row)
    $wsPrev.Hyperlinks.Add($cell, $urls[$i])
}

# The last source row (old "New" row 6 / new row 321) carries an empty-string
# (not blank) cadastreText cell; re-assert that so it keeps its text type.
$wsPrev.Range("E" + $lastRow).Value = "'"

# Adding a hyperlink / assigning .Value resets a cell's style to Excel's
# generic defaults; repaint the whole appended block with the formatting
# used by the rest of the table (copied from the row directly above it).
$lastRowBeforeBlock = $destStartRow - 1
$styleSrc = $wsPrev.Range("A" + $lastRowBeforeBlock + ":F" + $lastRowBeforeBlock)
$styleDestRange = $wsPrev.Range("A" + $destStartRow + ":F" + $lastRow)
$styleSrc.Copy()
$styleDestRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Remove the hyperlinks from "New" and clear its data rows, leaving only the header.
for ($i = 2; $i -le 6; $i++) {
    $cell = $wsNew.Range("A" + $i)
    if ($cell.Hyperlinks.Count() -gt 0) {
        $cell.Hyperlinks.Delete()
    }
}
$wsNew.Range("A2:F6").Delete()
